$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Start/End time values for the existing row 68
$ws.Range("B68").Value = 0.77569444444444446
$ws.Range("C68").Value = 0.96944444444444444

# Add the new daily power record for row 69 -- copy the formatting/formulas
# of an existing data row (67) down into the new row so the styles and
# relative formulas match what Excel would produce, then set its values.
$ws.Range("A67:F67").Copy()
$ws.Range("A69:F69").Insert(-4121)
$ws.Range("A69").Value = 43393
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = 0

# Row 70 only carries the calculated-column formulas (no Date/Start/End
# entered yet) -- again clone row 67's formatting/formulas, then clear the
# Date/Start/End cells so only the calculated columns remain populated.
$ws.Range("A67:F67").Copy()
$ws.Range("A70:F70").Insert(-4121)
$ws.Range("A70:C70").Clear()

# Grow the table to cover the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F70"))

# Update the saved selection state to match
$ws.Range("A70").Select()
